# Apply corrected Diebold-Mariano values to the "P_valores" and
# "Estadisticos_DM" sheets (Correcion a Diebold Mariano y revision de Cap1).

$wb = $excel.ActiveWorkbook

# --- Sheet "P_valores" ---------------------------------------------------
$wsP = $wb.Worksheets.Item("P_valores")

$pValores = @(
    @(1,                    0.3412959319354949,  0.5170447553144439,  0.3730469830511152,  0.9382657485717207),
    @(0.3412959319354949,   1,                    0.9269216456849638,  0.6099374695588731,  0.1495124185556427),
    @(0.5170447553144439,   0.9269216456849638,   1,                   0.8403870580241333,  0.2814721920435808),
    @(0.3730469830511152,   0.6099374695588731,   0.8403870580241333,  1,                    0.2430225132780002),
    @(0.9382657485717207,   0.1495124185556427,   0.2814721920435808,  0.2430225132780002,  1)
)

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt 5; $j++) {
        $col = $j + 2
        $wsP.Cells.Item($row, $col).Value = $pValores[$i][$j]
    }
}

# --- Sheet "Estadisticos_DM" ---------------------------------------------
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$estadisticos = @(
    @(0,                     0.9850970677827603,  0.6646961987255042,   0.9202110934151185,  -0.07885273352504625),
    @(-0.9850970677827603,   0,                    -0.09338414755428226, -0.5218426137174194, -1.525062306412314),
    @(-0.6646961987255042,   0.09338414755428226,  0,                    -0.2051770616041748, -1.120214612877305),
    @(-0.9202110934151185,   0.5218426137174194,   0.2051770616041748,   0,                    -1.218900349771393),
    @(0.07885273352504625,   1.525062306412314,    1.120214612877305,    1.218900349771393,   0)
)

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt 5; $j++) {
        $col = $j + 2
        $wsE.Cells.Item($row, $col).Value = $estadisticos[$i][$j]
    }
}

$wb.Save()
